$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values in B2:B5 (rounded values)
$ws.Range("B2").Value = 0.15
$ws.Range("B3").Value = 0.47
$ws.Range("B4").Value = 0.13
$ws.Range("B5").Value = 0.1

# Row 6 now becomes "Surprise" with new values, replacing "Disgust"
$ws.Range("A6").Value = "Surprise"
$ws.Range("B6").Value = 0.14
$ws.Range("C6").Value = 0.1

# Old row 7 (Surprise/0.128/0.5) is removed entirely - delete the whole row
$ws.Rows("7").Delete()
